# Generate Report for Handoff
# Adds a new row for file "e39f574c-b915-4285-95c4-dfdd9da38f93.md" to the
# Overview, zh-cn and de-de worksheets, expanding each table by one row.

$wb = $excel.ActiveWorkbook

$guidFile   = "e39f574c-b915-4285-95c4-dfdd9da38f93.md"
$guidPath   = "e2e\e39f574c-b915-4285-95c4-dfdd9da38f93.md"
$hashZh     = "e39f574c-b915-4285-95c4-dfdd9da38f93.6f9994cd4a132697e5c945d3f283e9462aaed486.zh-cn.xlf"
$hashDe     = "e39f574c-b915-4285-95c4-dfdd9da38f93.6f9994cd4a132697e5c945d3f283e9462aaed486.de-de.xlf"
$hoDateZh   = "2016-09-07 03:06:54"
$hoDateDe   = "2016-09-07 03:07:11"
$overviewDate = "2016-09-07 03:07:11"
$hbDateNone = "0001-01-01 00:00:00"
$status     = "Ready for handoff"
$hlUnderline = 2
$hlColor     = 15570276

function Set-HyperlinkFont($range) {
    $range.Font.Underline = $hlUnderline
    $range.Font.Color = $hlColor
}

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A9").Value = $guidFile
$wsOverview.Range("B9").Value = $guidPath
$wsOverview.Range("C9").Value = ".md"
$wsOverview.Range("D9").Value = ""
$wsOverview.Range("E9").Value = $status
$wsOverview.Range("F9").Value = $status
$wsOverview.Range("G9").Value = $overviewDate
$wsOverview.Range("G9").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("B9"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4cdf0d897301d1d31b120f7e4a52fd57a1cba7d2/e2e/e39f574c-b915-4285-95c4-dfdd9da38f93.md",
    "",
    "",
    $guidPath
) | Out-Null

$wsOverview.ListObjects.Item(1).Resize($wsOverview.Range("A1:G9"))

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A9").Value = $guidFile
$wsZh.Range("B9").Value = ".md"
$wsZh.Range("C9").Value = $status
$wsZh.Range("D9").Value = "e2e"
$wsZh.Range("E9").Value = "ht"
$wsZh.Range("F9").Value = "False"
$wsZh.Range("G9").Value = $hashZh
$wsZh.Range("H9").Value = $hoDateZh
$wsZh.Range("H9").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("I9").Value = ""
$wsZh.Range("J9").Value = ""
$wsZh.Range("K9").Value = $hbDateNone
$wsZh.Range("K9").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("L9").Value = ""
$wsZh.Range("M9").Value = "True"
$wsZh.Range("N9").Value = ""
$wsZh.Range("O9").Value = "False"
$wsZh.Range("P9").Value = ""

$wsZh.Hyperlinks.Add(
    $wsZh.Range("A9"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4cdf0d897301d1d31b120f7e4a52fd57a1cba7d2/e2e/e39f574c-b915-4285-95c4-dfdd9da38f93.md",
    "",
    "",
    $guidFile
) | Out-Null

$wsZh.ListObjects.Item(1).Resize($wsZh.Range("A1:P9"))

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A9").Value = $guidFile
$wsDe.Range("B9").Value = ".md"
$wsDe.Range("C9").Value = $status
$wsDe.Range("D9").Value = "e2e"
$wsDe.Range("E9").Value = "ht"
$wsDe.Range("F9").Value = "False"
$wsDe.Range("G9").Value = $hashDe
$wsDe.Range("H9").Value = $hoDateDe
$wsDe.Range("H9").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("I9").Value = ""
$wsDe.Range("J9").Value = ""
$wsDe.Range("K9").Value = $hbDateNone
$wsDe.Range("K9").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("L9").Value = ""
$wsDe.Range("M9").Value = "True"
$wsDe.Range("N9").Value = ""
$wsDe.Range("O9").Value = "False"
$wsDe.Range("P9").Value = ""

$wsDe.Hyperlinks.Add(
    $wsDe.Range("A9"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4cdf0d897301d1d31b120f7e4a52fd57a1cba7d2/e2e/e39f574c-b915-4285-95c4-dfdd9da38f93.md",
    "",
    "",
    $guidFile
) | Out-Null

$wsDe.ListObjects.Item(1).Resize($wsDe.Range("A1:P9"))

Write-Host "Added handoff row for $guidFile to Overview, zh-cn and de-de sheets."
